{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst p = paragraphs.items[2];\nconst range = p.getRange(\"Start\");\nreturn \"got range\";\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs.Item(1)\n$p.Format.WidowControl = $p.Format.WidowControl\nWrite-Output \"noop set done\"\n"}
